$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New recipe rows to append to the Recipes_VALLEJO table (columns A..L)
$data = @(
    @("'300059454", "DORITOS 3D POPMIX HBSF 120GRX12X1", "No", "No", "No", "No", "NA", "NA", "NA", "No", "No", "Packaging Pellet"),
    @("'300059455", "DORITOS 3D POPMIX HBSF RTD 120GRX12X1", "No", "No", "No", "No", "NA", "NA", "NA", "No", "No", "Packaging Pellet"),
    @("'300059450", "DORITOS 3D POPMIX HBSF RTDBAU 120GRX12X1", "No", "No", "No", "No", "NA", "NA", "NA", "No", "No", "Packaging Pellet")
)

$startRow = 2
$lastCol = 12

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $row[$c]
    }
}

# Copy header row formatting (style) onto the newly added rows so they
# match the look of the existing data/header row.
$headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $lastCol))
$headerRange.Copy()
$targetRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($startRow + $data.Count - 1, $lastCol))
$targetRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Grow the table (ListObject) to cover the newly populated rows.
$table = $ws.ListObjects.Item(1)
$newRef = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($startRow + $data.Count - 1, $lastCol))
$table.Resize($newRef)
